$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row 18 data
$ws.Range("B18").Value = (Get-Date -Year 2024 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = "Tein asetukset ja aloin tyylitelemään."

# Update the sum formula to include the new row
$ws.Range("C19").Formula = "=SUM(C6:C18)"

# Move the selection to D18 to match the saved selection state
$ws.Range("D18").Select()
